# "Generate Report for Handback"
#
# A new handback was processed for 056f69cc-0802-4d88-954c-56dff0abc57f
# (row 7) on both the zh-cn and de-de target-language sheets. The
# handback's file version turned out to be stale, so the row records the
# "Latest Target File" / "Latest Handback File" / "Latest Handback
# DateTime" / "Error Detail" columns accordingly, and the target-file
# cell (column I) becomes a hyperlink into the target-language repo,
# exactly like column A already is.

$wb = $excel.ActiveWorkbook

$handbackFile        = "056f69cc-0802-4d88-954c-56dff0abc57f.md"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/946791313f895cde34c0973ea5aae77c328c906a/e2e/056f69cc-0802-4d88-954c-56dff0abc57f.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f17636ca5105280361fcc6995025644c200a9a9c/e2e/056f69cc-0802-4d88-954c-56dff0abc57f.md."

function Update-HandbackRow {
    param(
        [string]$SheetName,
        [string]$OrgRepoSuffix,
        [string]$HandbackDateTime
    )

    $ws = $wb.Worksheets.Item($SheetName)
    Write-Host "Updating sheet $SheetName ($OrgRepoSuffix) @ $HandbackDateTime"

    # Column G = "Latest Handoff File" already holds the xlf file name for
    # this row; column J = "Latest Handback File" reuses the same text.
    $xliffName = $ws.Range("G7").Value

    # Column I = "Latest Target File" - becomes a hyperlink, same display
    # text/target as the A7 handoff-markdown hyperlink but pointed at the
    # language-specific repo.
    $ws.Range("I7").Value = $handbackFile
    $targetUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0-$OrgRepoSuffix/blob/946791313f895cde34c0973ea5aae77c328c906a/e2e/$handbackFile"
    $ws.Hyperlinks.Add($ws.Range("I7"), $targetUrl, "", "", $handbackFile)

    # Column J = "Latest Handback File"
    $ws.Range("J7").Value = $xliffName

    # Column K = "Latest Handback DateTime"
    $ws.Range("K7").Value = $HandbackDateTime

    # Column P = "Error Detail"
    $ws.Range("P7").Value = $errorDetail
}

Update-HandbackRow "zh-cn" "zhcn" "2016-08-26 12:55:55"
Update-HandbackRow "de-de" "dede" "2016-08-26 12:56:08"
